# The sheet gains one new data row: a row is inserted at row 756 (pushing
# the former rows 756-838 down to 757-839, and growing the used range from
# A1:T838 to A1:T839). The new row carries the same "Macroferia Regional de
# Talca" / "Frutilla" category metadata as its neighbours, with its own
# date, quality, volume, price and derived $/Kg figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 756; Excel shifts rows 756:838 down to 757:839.
$ws.Rows(756).Insert()

# Populate the newly inserted row 756 with its data.
$ws.Range("A756").Value = 5
$ws.Range("B756").Value = "Macroferia Regional de Talca"
$ws.Range("C756").Value = "Maule"
$ws.Range("D756").Value = 45194
$ws.Range("E756").Value = 7
$ws.Range("F756").Value = "Fruta"
$ws.Range("G756").Value = 100101
$ws.Range("H756").Value = "Berries"
$ws.Range("I756").Value = 100112025
$ws.Range("J756").Value = "Frutilla"
$ws.Range("K756").Value = "Sin especificar"
$ws.Range("L756").Value = "Segunda"
$ws.Range("M756").Value = 30
$ws.Range("N756").Value = 12000
$ws.Range("O756").Value = 12000
$ws.Range("P756").Value = 12000
$ws.Range("Q756").Value = "`$/bandeja 7 kilos"
$ws.Range("R756").Value = "Provincia de Melipilla"
$ws.Range("S756").Value = 1714
$ws.Range("T756").Value = 7
